$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.605.74"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "2.600.82"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "567.10"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.57%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.93"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.03%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "2.623.69"
$ws.Range("E9").Value = "  +1.54%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.54"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.371"
$c.ClearFormats()
$ws.Range("E12").Value = "  +7.76%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.153"
$c.ClearFormats()
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D14").Value = "3.062.23"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "60.627.62"
$ws.Range("E15").Value = "  +2.46%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "23.40"
$c.ClearFormats()
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "2.613.14"
$ws.Range("E18").Value = "  +1.34%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.19"
$c.ClearFormats()
$ws.Range("E19").Value = "  +8.39%  "
$ws.Range("E20").Value = "  +2.53%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "347.41"
$c.ClearFormats()
$ws.Range("E21").Value = "  +3.07%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.04"
$c.ClearFormats()
$ws.Range("E22").Value = "  +11.13%  "
$ws.Range("E23").Value = "  -0.04%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.532"
$c.ClearFormats()
$ws.Range("E24").Value = "  +16.16%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "63.51"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.80%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.41%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.ClearFormats()
$ws.Range("E27").Value = "  -1.64%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.ClearFormats()
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("E29").Value = "  +1.79%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.ClearFormats()
$ws.Range("E30").Value = "  +8.37%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +4.37%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "161.62"
$c.ClearFormats()
$ws.Range("E33").Value = "  +0.33%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "19.51"
$c.ClearFormats()
$ws.Range("E34").Value = "  +3.12%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.25"
$c.ClearFormats()
$ws.Range("E35").Value = "  +5.82%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.970"
$c.ClearFormats()
$ws.Range("E36").Value = "  +11.38%  "
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("E38").Value = "  +7.71%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "37.84"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.13%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.ClearFormats()
$ws.Range("E40").Value = "  +4.42%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.848"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.96%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "295.91"
$c.ClearFormats()
$ws.Range("E42").Value = "  +0.56%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "140.11"
$c.ClearFormats()
$ws.Range("E43").Value = "  +6.54%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  +2.19%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0979"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.66%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0547"
$c.ClearFormats()
$ws.Range("E47").Value = "  +2.26%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "19.55"
$c.ClearFormats()
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("E50").Value = "  +0.59%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.87"
$c.ClearFormats()
$ws.Range("E51").Value = "  +8.26%  "
